# updated CB_API and Dash
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("S2").Value = 36.5625
$ws.Range("T2").Value = 43.0078125
$ws.Range("U2").Value = 48.07125
$ws.Range("V2").Value = 38.287109375

# Row 14 updates
$ws.Range("J14").Value = 40.50396634615385
$ws.Range("K14").Value = 42.84707661290322
$ws.Range("L14").Value = 50.5265625
$ws.Range("M14").Value = 54.85258928571428
$ws.Range("N14").Value = 36.14714673913043
$ws.Range("O14").Value = 38.66835937499999
$ws.Range("P14").Value = 39.88155241935483
$ws.Range("Q14").Value = 38.23306451612903
$ws.Range("R14").Value = 41.12668269230769
$ws.Range("S14").Value = 49.618125
$ws.Range("T14").Value = 54.7303125
$ws.Range("U14").Value = 44.99442567567567
$ws.Range("V14").Value = 42.63515624999999
$ws.Range("W14").Value = 36.54375
